# CucumberTest/src/TestData/TestData.xlsx
# "Added steps in feature file and basic asserts"
#
# Adds a new "Driver_Type" column (with "Occasional" values) to the
# CustomerData sheet, and tidies up a few sheet selections / the
# ReinstatementData page orientation left over from the author's Excel
# session.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# VehicleData: selection moved from A1:C3 to G1:G3 (active cell G3)
# ---------------------------------------------------------------------
$wsVehicle = $wb.Worksheets.Item(3)
$wsVehicle.Activate()
$wsVehicle.Range("G1:G3").Select()

# ---------------------------------------------------------------------
# EndorsementData: selection moved from C7 to A3 (and it stops being the
# sheet that is on top when the workbook is reopened)
# ---------------------------------------------------------------------
$wsEndorsement = $wb.Worksheets.Item(4)
$wsEndorsement.Activate()
$wsEndorsement.Range("A3").Select()

# ---------------------------------------------------------------------
# ReinstatementData: selection moved from D8 to E1, and the sheet is set
# to print in portrait orientation
# ---------------------------------------------------------------------
$wsReinstatement = $wb.Worksheets.Item(7)
$wsReinstatement.Activate()
$wsReinstatement.Range("E1").Select()
$wsReinstatement.PageSetup.Orientation = 1

# ---------------------------------------------------------------------
# CustomerData: new "Driver_Type" column (M) with "Occasional" for both
# data rows; this becomes the active sheet, selection on the new header
# cell M1
# ---------------------------------------------------------------------
$wsCustomer = $wb.Worksheets.Item(2)
$wsCustomer.Activate()

$wsCustomer.Range("M1").Value = "Driver_Type"
$wsCustomer.Range("M2").Value = "Occasional"
$wsCustomer.Range("M3").Value = "Occasional"
$wsCustomer.Columns.Item(13).ColumnWidth = 12

$wsCustomer.Range("M1").Select()
